# Weekly price-data update for "Hortaliza, Terminal Hortofrutícola Agro
# Chillán - Acelga": insert two new observation rows (week of 2023-05-29 /
# Excel serial 45075) ahead of the existing data, pushing the previous rows
# 447-462 down to 449-464.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 447:448 (existing rows shift down to 449:450, etc.)
$ws.Range("A447:A448").EntireRow.Insert()

# --- New row 447 ("Primera") ---
$ws.Cells.Item(447, 1).Value = 7
$ws.Cells.Item(447, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(447, 3).Value = "Ñuble"
$ws.Cells.Item(447, 4).Value = 45075
$ws.Cells.Item(447, 5).Value = 16
$ws.Cells.Item(447, 6).Value = 100112009
$ws.Cells.Item(447, 7).Value = "Acelga"
$ws.Cells.Item(447, 8).Value = "Sin especificar"
$ws.Cells.Item(447, 9).Value = "Primera"
$ws.Cells.Item(447, 10).Value = 100
$ws.Cells.Item(447, 11).Value = 700
$ws.Cells.Item(447, 12).Value = 700
$ws.Cells.Item(447, 13).Value = 700
$ws.Cells.Item(447, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(447, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(447, 16).Value = 700
$ws.Cells.Item(447, 17).Value = 1
$ws.Cells.Item(447, 18).Value = "Hortaliza"

# --- New row 448 ("Segunda") ---
$ws.Cells.Item(448, 1).Value = 7
$ws.Cells.Item(448, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(448, 3).Value = "Ñuble"
$ws.Cells.Item(448, 4).Value = 45075
$ws.Cells.Item(448, 5).Value = 16
$ws.Cells.Item(448, 6).Value = 100112009
$ws.Cells.Item(448, 7).Value = "Acelga"
$ws.Cells.Item(448, 8).Value = "Sin especificar"
$ws.Cells.Item(448, 9).Value = "Segunda"
$ws.Cells.Item(448, 10).Value = 150
$ws.Cells.Item(448, 11).Value = 500
$ws.Cells.Item(448, 12).Value = 500
$ws.Cells.Item(448, 13).Value = 500
$ws.Cells.Item(448, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(448, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(448, 16).Value = 500
$ws.Cells.Item(448, 17).Value = 1
$ws.Cells.Item(448, 18).Value = "Hortaliza"
